$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 415; existing rows 415-448 shift down to 416-449.
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row 415 with the new weekly data record.
$ws.Cells.Item(415, 1).Value = 6
$ws.Cells.Item(415, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(415, 3).Value = "Metropolitana"
$ws.Cells.Item(415, 4).Value = 44826
$ws.Cells.Item(415, 5).Value = 13
$ws.Cells.Item(415, 6).Value = 100112032
$ws.Cells.Item(415, 7).Value = "Zapallo italiano"
$ws.Cells.Item(415, 8).Value = "Sin especificar"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 720
$ws.Cells.Item(415, 11).Value = 15000
$ws.Cells.Item(415, 12).Value = 16000
$ws.Cells.Item(415, 13).Value = 15681
$ws.Cells.Item(415, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(415, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(415, 16).Value = 314
$ws.Cells.Item(415, 17).Value = 50
$ws.Cells.Item(415, 18).Value = "Hortaliza"

"row 415 inserted and populated"
